# "save data done + era data updated"
# Add a new "Save" column (H) to the sheet, matching the header styling
# already used by B1:G1, and populate the data row with the saved flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", formatted like the rest of the header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# New data cell H2 = 1 (save flag for this row).
$ws.Range("H2").Value = 1
